# Update column F (dSF) values for specific rows per repulled data / mean calculation fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 3
    6  = 2
    8  = 3
    9  = 0
    13 = 2
    14 = -1
    17 = 0
    19 = 1
    22 = 0
    29 = -4
    36 = 1
    37 = -3
    38 = 3
    39 = 3
    40 = -3
    48 = 2
    51 = 3
    59 = 0
    64 = -1
    66 = 6
    71 = -7
    72 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
